# Applies the 202212 report-data refresh described in the commit:
#   - "Printed on" timestamp bumped
#   - chart callout dollar figures bumped by $1-2 (rounding refresh)
#   - underlying computed ratios/balances on Sheet1/2/4/5 refreshed
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet1 (visible) - recomputed ratio / balance figures
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("N25").Value = 0.11018268964095
$ws1.Range("N26").Value = 0.00295119505621494
$ws1.Range("N27").Value = 0.00462235370250532
$ws1.Range("N28").Value = 0.000426678803308199
$ws1.Range("O46").Value = 0.0510444982208855
$ws1.Range("P46").Value = 0.0340174276673517
$ws1.Range("Q46").Value = 0.0258439252916966
$ws1.Range("R46").Value = 0.0257595709597618
$ws1.Range("N47").Value = 0.0350540215252864
$ws1.Range("N48").Value = -0.00605739821251241
$ws1.Range("N49").Value = -0.309856122576091
$ws1.Range("N51").Value = -0.0627616119100511
$ws1.Range("F55").Value = 13228.742366378
$ws1.Range("G55").Value = 13568.598339558801
$ws1.Range("H55").Value = 13858.3446171692
$ws1.Range("I55").Value = 14100.2295742386
$ws1.Range("J55").Value = 13688.978724336199
$ws1.Range("O55").Value = 24690.073260616
$ws1.Range("E56").Value = 106497.017
$ws1.Range("F56").Value = 106497.017
$ws1.Range("G56").Value = 106497.017
$ws1.Range("H56").Value = 106497.017
$ws1.Range("I56").Value = 106497.017
$ws1.Range("J56").Value = 106497.017
$ws1.Range("O56").Value = 4.66836062820837
$ws1.Range("E57").Value = 120242.999
$ws1.Range("F57").Value = 120242.999
$ws1.Range("G57").Value = 120242.999
$ws1.Range("H57").Value = 120242.999
$ws1.Range("I57").Value = 120242.999
$ws1.Range("J57").Value = 120242.999
$ws1.Range("O57").Value = 0.17543308486614
$ws1.Range("E58").Value = 140622.752
$ws1.Range("F58").Value = 140751.494366378
$ws1.Range("G58").Value = 141091.350339559
$ws1.Range("H58").Value = 141381.09661717
$ws1.Range("I58").Value = 141622.98157423898
$ws1.Range("J58").Value = 141211.730724336
$ws1.Range("O58").Value = 0.0691119416445477
$ws1.Range("E59").Value = 136217.752
$ws1.Range("F59").Value = 136346.494366378
$ws1.Range("G59").Value = 136686.350339559
$ws1.Range("H59").Value = 136976.09661717
$ws1.Range("I59").Value = 137217.98157423898
$ws1.Range("J59").Value = 136806.730724336
$ws1.Range("O59").Value = 0.0431157477391154
$ws1.Range("F60").Value = 16378.742366378
$ws1.Range("G60").Value = 16718.5983395587
$ws1.Range("H60").Value = 17008.3446171697
$ws1.Range("I60").Value = 17250.2295742387
$ws1.Range("J60").Value = 16838.9787243363
$ws1.Range("O60").Value = -0.0500038327141156
$ws1.Range("F61").Value = 1692.1105284939101
$ws1.Range("G61").Value = 1693.8134130987398
$ws1.Range("H61").Value = 1704.36675963478
$ws1.Range("I61").Value = 1703.83545298323
$ws1.Range("J61").Value = 6794.1261542106595
$ws1.Range("O61").Value = -0.104776741352632
$ws1.Range("F62").Value = 187.881097147783
$ws1.Range("G62").Value = 262.639883643089
$ws1.Range("H62").Value = 337.956587920571
$ws1.Range("I62").Value = 401.867459941974
$ws1.Range("J62").Value = 1190.34502865342
$ws1.Range("O62").Value = -0.166440467163279
$ws1.Range("F63").Value = 1504.22943134612
$ws1.Range("G63").Value = 1431.17352945565
$ws1.Range("H63").Value = 1366.41017171421
$ws1.Range("I63").Value = 1301.96799304126
$ws1.Range("J63").Value = 5603.78112555724
$ws1.Range("O63").Value = -0.228751582644323
$ws1.Range("F64").Value = 361.74278163735704
$ws1.Range("G64").Value = 306.94014321042704
$ws1.Range("H64").Value = 258.358128810836
$ws1.Range("I64").Value = 210.017045806315
$ws1.Range("J64").Value = 1137.05809946494
$ws1.Range("O65").Value = -0.0831690682530268

# "Printed on" banner (R8, shared string)
$ws1.Range("R8").Value = "Printed on:2023-02-10 14:32"

# ---------------------------------------------------------------
# Sheet2 (hidden) - trend helper figures
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("B2").Value = 0.558629075077603
$ws2.Range("C2").Value = 0.873010916823471
$ws2.Range("B3").Value = 1.12990322877871
$ws2.Range("C3").Value = 0.996835584305405
$ws2.Range("B4").Value = 1.29556240919319
$ws2.Range("C4").Value = 1.09558697347087
$ws2.Range("B5").Value = 1.2875701627367
$ws2.Range("C5").Value = 1.51098840930809
$ws2.Range("B6").Value = 1.15416622238456
$ws2.Range("C6").Value = 0.997580077927072

# ---------------------------------------------------------------
# Sheet4 (hidden) - balance-sheet history row for 202212
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("B9").Value = 106497.017
$ws4.Range("C9").Value = 4460.018
$ws4.Range("D9").Value = 67968.002
$ws4.Range("E9").Value = 14827.999
$ws4.Range("G9").Value = 3072.002
$ws4.Range("K9").Value = 6159.994
$ws4.Range("N9").Value = 6522.998
$ws4.Range("Q9").Value = 3486.004

# ---------------------------------------------------------------
# Sheet5 (hidden) - single trend figure for 202212
# ---------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Sheet5")
$ws5.Range("B8").Value = 1518.017

# ---------------------------------------------------------------
# Chart callout text boxes anchored on Sheet1 (xl/drawings/drawing1.xml)
# Each shape holds "<pct>%\n$<amount>" - bump the dollar amount only.
# ---------------------------------------------------------------
$ws1.Shapes.Item("ReturnonAssets").TextFrame2.TextRange.Text = "0.80%`n`$1,137"
$ws1.Shapes.Item("ProvisionforTax").TextFrame2.TextRange.Text = "0.27%`n`$379"
$ws1.Shapes.Item("OperatingProfit").TextFrame2.TextRange.Text = "1.07%`n`$1,516"
$ws1.Shapes.Item("NetInterestMargin").TextFrame2.TextRange.Text = "3.97%`n`$5,604"
$ws1.Shapes.Item("InterestExpense").TextFrame2.TextRange.Text = "0.84%`n`$1,190"
$ws1.Shapes.Item("Deposits").TextFrame2.TextRange.Text = "0.78%`n`$1,099"
$ws1.Shapes.Item("Borrowings").TextFrame2.TextRange.Text = "0.06%`n`$91"
